$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the "From" value for rule R20 (row 10, column C) from 18 to 100
$ws.Range("C10").Value = 100
